# Update "horarios-141" workbook with the latest scrape results.
# Sheet 1: LP1912
# Sheet 2: LP1912-215
# Sheet 3: 6203-6173

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:35:49"
$ws1.Range("A3").Value = "Total filas: 15"

$sheet1Data = @(
    @("03:35:49", "03:38", "14_ABASTO",      3,   "LP1912"),
    @("02:58:51", "03:48", "14_ABASTO",      50,  "LP1912"),
    @("02:21:47", "03:56", "14_ABASTO",      95,  "LP1912"),
    @("03:35:49", "04:01", "81_EL PELIGRO",  26,  "LP1912"),
    @("03:35:49", "04:45", "215A_EL PATO",   70,  "LP1912"),
    @("03:35:49", "04:53", "11_ETCHEVERRY",  78,  "LP1912"),
    @("03:35:49", "05:16", "17_ROMERO",      101, "LP1912"),
    @("03:35:49", "05:22", "23_HERNANDEZ",   107, "LP1912"),
    @("03:35:49", "05:28", "14_ABASTO",      113, "LP1912"),
    @("03:35:49", "05:34", "215B_EL PATO",   119, "LP1912")
)

$row = 11
foreach ($r in $sheet1Data) {
    $ws1.Cells.Item($row, 1).Value = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Cells.Item($row, 4).Value = $r[3]
    $ws1.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:35:49"
$ws2.Range("A3").Value = "Total filas: 6"

$ws2.Cells.Item(10, 1).Value = "03:35:49"
$ws2.Cells.Item(10, 2).Value = "04:45"
$ws2.Cells.Item(10, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(10, 4).Value = 70
$ws2.Cells.Item(10, 5).Value = "LP1912"

$ws2.Cells.Item(11, 1).Value = "03:35:49"
$ws2.Cells.Item(11, 2).Value = "05:34"
$ws2.Cells.Item(11, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(11, 4).Value = 119
$ws2.Cells.Item(11, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:35:49"
